$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.8179
$ws.Range("A9").Value = -21.96530000000001
$ws.Range("B9").Value = 5.137000000000006
$ws.Range("D9").Value = -8.446699999999996
$ws.Range("A18").Value = -22.224
$ws.Range("A20").Value = -20.80629999999998
$ws.Range("B23").Value = 8.950899999999997
$ws.Range("B24").Value = 5.118600000000006
$ws.Range("B26").Value = 4.874200000000003
$ws.Range("A27").Value = -21.9059
$ws.Range("D32").Value = -7.243599999999997
$ws.Range("B34").Value = 9.417700000000005
$ws.Range("B35").Value = 8.780100000000008
$ws.Range("D38").Value = -7.462399999999998
$ws.Range("D45").Value = -7.066699999999998
$ws.Range("B48").Value = 5.566500000000004
$ws.Range("D51").Value = -8.2704
$ws.Range("B52").Value = 5.672599999999997
$ws.Range("D57").Value = -7.879999999999995
$ws.Range("D64").Value = -7.337699999999992
$ws.Range("B66").Value = 5.449699999999996
$ws.Range("B67").Value = 5.471799999999997
$ws.Range("A69").Value = -21.65639999999998
$ws.Range("A76").Value = -19.55989999999997
$ws.Range("B80").Value = 9.604199999999993
$ws.Range("A82").Value = -21.63830000000002
$ws.Range("D93").Value = -6.922199999999993
$ws.Range("B99").Value = 6.248600000000003
